# BAU Rng Anxiety and Charge Time Shadow Costs.xlsx
# "updated files from canada" - re-order the BRAaCTSC year columns (B:D) so
# the years run chronologically (2018, 2019, 2020 instead of 2020, 2018,
# 2019), fix up the TREND() formulas that ride along with each year column,
# unhide the (previously hidden) helper column B, and move the active
# selection from F7 to D7 - matching the edits captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BRAaCTSC")

# --- Row 1: re-order the year headers in B1:D1 -----------------------------
# Before: B1=2020, C1=2018, D1=2019
# After : B1=2018, C1=2019, D1=2020
$ws.Range("B1").Value = 2018
$ws.Range("C1").Value = 2019
$ws.Range("D1").Value = 2020

# --- Row 2: the TREND() formulas shift together with the headers ----------
# Each formula keeps referencing "the next few columns" / "this column's
# header", so after the re-order the formula text itself shifts by one
# column (wrapping from D back to B).
$ws.Range("B2").Formula = "=TREND(C2:H2,C1:H1,B1)"
$ws.Range("C2").Formula = "=TREND(D2:I2,D1:I1,C1)"
$ws.Range("D2").Formula = "=TREND(Calculations!`$B`$33:`$B`$34,Calculations!`$A`$33:`$A`$34,BRAaCTSC!D1)"

# --- Unhide column B (was width 0 / hidden) --------------------------------
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(2).ColumnWidth = 8

# --- Update the active selection on the BRAaCTSC sheet: F7 -> D7 ----------
$ws.Activate() | Out-Null
$ws.Range("D7").Select() | Out-Null
